# Week7 slide: change the dice-rolling example so that a single 20-sided
# die is rolled 5 times instead of 10 times, and the randint() upper bound
# used for that die becomes 20 instead of 5.
#
#   for i in range(10):              ->   for i in range(5):
#       results.append(random.randint(1,5))  ->  ...randint(1,20))

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(21)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 19 (1-based): "for i in range(10):" -- third run holds " in range(10):"
$forPara = $tr.Paragraphs(19, 1)
$forPara.Runs(3, 1).Text = " in range(5):"

# Paragraph 20 (1-based): "    results.append(random.randint(1,5))"
# -- fifth run holds "(1,5))"
$appendPara = $tr.Paragraphs(20, 1)
$appendPara.Runs(5, 1).Text = "(1,20))"
